# The sheet holds per-innings stats for Dinesh Karthik (Kolkata Knight
# Riders). The runs/balls/fours/sixes columns (C:F) for rows 2-8 got
# reshuffled - this reproduces the new values cell by cell.
#
# All of these columns are stored as text (numbers-as-text, e.g. "21"
# rather than 21), so each value is entered with a leading apostrophe to
# force a text entry instead of letting Excel infer a numeric type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Formula = "'0"
$ws.Range("D2").Formula = "'2"
$ws.Range("E2").Formula = "'0"

$ws.Range("C3").Formula = "'29"
$ws.Range("D3").Formula = "'14"
$ws.Range("E3").Formula = "'2"
$ws.Range("F3").Formula = "'2"

$ws.Range("D4").Formula = "'14"
$ws.Range("E4").Formula = "'0"

$ws.Range("C5").Formula = "'21"
$ws.Range("D5").Formula = "'10"
$ws.Range("E5").Formula = "'3"

$ws.Range("C6").Formula = "'4"
$ws.Range("D6").Formula = "'8"
$ws.Range("E6").Formula = "'1"

$ws.Range("C7").Formula = "'0"
$ws.Range("D7").Formula = "'1"
$ws.Range("E7").Formula = "'0"
$ws.Range("F7").Formula = "'0"

$ws.Range("C8").Formula = "'3"
$ws.Range("D8").Formula = "'6"
